# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
# Mapping of cell -> new value (same update applied identically on both sheets).
$updates = @{
    "F2"  = 248
    "F4"  = 285
    "F6"  = 284
    "F7"  = 6813
    "F8"  = 61
    "F13" = 15
    "F16" = 232
    "F17" = 589
    "F18" = 69
}

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
